# Auto-generated Excel COM-interop script applying numeric updates
# to the Hyperion_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5038.016
$ws.Range("I17").Value = 800
$ws.Range("K17").Value = 2400
$ws.Range("M17").Value = -2232
$ws.Range("H18").Value = 4810929
$ws.Range("I18").Value = 7695252.5
$ws.Range("J18").Value = 3723
$ws.Range("K18").Value = 7695252.5
$ws.Range("L18").Value = 3723
$ws.Range("M18").Value = -7694968.5
$ws.Range("N18").Value = -4291
$ws.Range("H33").Value = 1236.125
$ws.Range("I33").Value = 1177.8
$ws.Range("J33").Value = 1333.3334
$ws.Range("K33").Value = 1177.8
$ws.Range("L33").Value = 1333.3334
$ws.Range("M33").Value = -948.8
$ws.Range("N33").Value = -1791.3334
$ws.Range("H86").Value = 2279.3
$ws.Range("I86").Value = 2724.25
$ws.Range("J86").Value = 1982.6666
$ws.Range("K86").Value = 2724.25
$ws.Range("L86").Value = 1982.6666
$ws.Range("M86").Value = -1601.25
$ws.Range("N86").Value = -4228.6666
$ws.Range("H89").Value = 2279.3
$ws.Range("I89").Value = 2724.25
$ws.Range("J89").Value = 1982.6666
$ws.Range("K89").Value = 13621.25
$ws.Range("L89").Value = 9913.333000000001
$ws.Range("M89").Value = -8005.25
$ws.Range("N89").Value = -21145.333
$ws.Range("H98").Value = 2794.3333
$ws.Range("I98").Value = 2619.4827
$ws.Range("J98").Value = 3518.7144
$ws.Range("K98").Value = 2619.4827
$ws.Range("L98").Value = 3518.7144
$ws.Range("M98").Value = -1121.4827
$ws.Range("N98").Value = -6514.7144
$ws.Range("H122").Value = 2794.3333
$ws.Range("I122").Value = 2619.4827
$ws.Range("J122").Value = 3518.7144
$ws.Range("K122").Value = 7858.4481
$ws.Range("L122").Value = 10556.1432
$ws.Range("M122").Value = -5408.4481
$ws.Range("N122").Value = -15456.1432
$ws.Range("H137").Value = 56866.383
$ws.Range("I137").Value = 78332.164
$ws.Range("J137").Value = 5348.5
$ws.Range("K137").Value = 234996.492
$ws.Range("L137").Value = 16045.5
$ws.Range("M137").Value = -232446.492
$ws.Range("N137").Value = -21145.5
$ws.Range("H138").Value = 4398.0557
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4398.0557
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 13194.1671
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -23474.1671
$ws.Range("H141").Value = 7569.7827
$ws.Range("I141").Value = 3709.1714
$ws.Range("K141").Value = 11127.5142
$ws.Range("M141").Value = -5947.514200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6258147
$ws.Range("I45").Value = 10276950
$ws.Range("K45").Value = 10276950
$ws.Range("M45").Value = -10276573
$ws.Range("H61").Value = 2692.6316
$ws.Range("I61").Value = 2632.6
$ws.Range("J61").Value = 2808.077
$ws.Range("K61").Value = 2632.6
$ws.Range("L61").Value = 2808.077
$ws.Range("M61").Value = -2420.6
$ws.Range("N61").Value = -3232.077
$ws.Range("H123").Value = 74644
$ws.Range("J123").Value = 74644
$ws.Range("L123").Value = 74644
$ws.Range("N123").Value = -84444
$ws.Range("H132").Value = 2647.4927
$ws.Range("I132").Value = 2040.425
$ws.Range("K132").Value = 6121.275
$ws.Range("M132").Value = -3591.275
$ws.Range("H136").Value = 2692.6316
$ws.Range("I136").Value = 2632.6
$ws.Range("J136").Value = 2808.077
$ws.Range("K136").Value = 7897.799999999999
$ws.Range("L136").Value = 8424.231
$ws.Range("M136").Value = -5347.799999999999
$ws.Range("N136").Value = -13524.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 54070
$ws.Range("J68").Value = 54070
$ws.Range("L68").Value = 54070
$ws.Range("N68").Value = -55692
$ws.Range("H71").Value = 54070
$ws.Range("J71").Value = 54070
$ws.Range("L71").Value = 162210
$ws.Range("N71").Value = -170322
$ws.Range("H86").Value = 3709970
$ws.Range("I86").Value = 4354634.5
$ws.Range("K86").Value = 4354634.5
$ws.Range("M86").Value = -4353511.5
$ws.Range("H89").Value = 3709970
$ws.Range("I89").Value = 4354634.5
$ws.Range("K89").Value = 21773172.5
$ws.Range("M89").Value = -21767556.5
$ws.Range("H107").Value = 4203854
$ws.Range("I107").Value = 4763805.5
$ws.Range("K107").Value = 4763805.5
$ws.Range("M107").Value = -4761885.5
$ws.Range("H134").Value = 2956.2031
$ws.Range("I134").Value = 1297.4286
$ws.Range("J134").Value = 6122.9546
$ws.Range("K134").Value = 3892.2858
$ws.Range("L134").Value = 18368.8638
$ws.Range("M134").Value = -1357.2858
$ws.Range("N134").Value = -23438.8638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 538.9091
$ws.Range("I22").Value = 347
$ws.Range("K22").Value = 347
$ws.Range("M22").Value = 3
$ws.Range("H31").Value = 18767.072
$ws.Range("I31").Value = 1744.591
$ws.Range("J31").Value = 29781.617
$ws.Range("K31").Value = 1744.591
$ws.Range("L31").Value = 29781.617
$ws.Range("M31").Value = -1449.591
$ws.Range("N31").Value = -30371.617
$ws.Range("H34").Value = 18767.072
$ws.Range("I34").Value = 1744.591
$ws.Range("J34").Value = 29781.617
$ws.Range("K34").Value = 1744.591
$ws.Range("L34").Value = 29781.617
$ws.Range("M34").Value = -1542.591
$ws.Range("N34").Value = -30185.617
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H132").Value = 23741.047
$ws.Range("I132").Value = 1945.9584
$ws.Range("J132").Value = 52801.168
$ws.Range("K132").Value = 5837.8752
$ws.Range("L132").Value = 158403.504
$ws.Range("M132").Value = -3307.8752
$ws.Range("N132").Value = -163463.504

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 50005600
$ws.Range("I56").Value = 50005600
$ws.Range("K56").Value = 50005600
$ws.Range("M56").Value = -50005070
$ws.Range("H81").Value = 7818.8667
$ws.Range("I81").Value = 937.6667
$ws.Range("K81").Value = 2813.0001
$ws.Range("M81").Value = -1690.0001
$ws.Range("H84").Value = 7818.8667
$ws.Range("I84").Value = 937.6667
$ws.Range("K84").Value = 8439.0003
$ws.Range("M84").Value = -2823.0003
$ws.Range("H92").Value = 1044.8572
$ws.Range("I92").Value = 1044.8572
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 3134.5716
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1886.5716
$ws.Range("N92").ClearContents()
$ws.Range("H121").Value = 23810116
$ws.Range("I121").Value = 41666864
$ws.Range("K121").Value = 125000592
$ws.Range("M121").Value = -124999282
$ws.Range("H129").Value = 1230.2858
$ws.Range("I129").Value = 1065.909
$ws.Range("K129").Value = 3197.727
$ws.Range("M129").Value = 1802.273
$ws.Range("H131").Value = 8015389
$ws.Range("J131").Value = 8133120
$ws.Range("L131").Value = 24399360
$ws.Range("N131").Value = -24409440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 19180.25
$ws.Range("I58").Value = 12814.667
$ws.Range("J58").Value = 22999.6
$ws.Range("K58").Value = 12814.667
$ws.Range("L58").Value = 22999.6
$ws.Range("M58").Value = -12537.667
$ws.Range("N58").Value = -23553.6
$ws.Range("H97").Value = 917054.0600000001
$ws.Range("I97").Value = 1036469.8
$ws.Range("K97").Value = 1036469.8
$ws.Range("M97").Value = -1035973.8
$ws.Range("H107").Value = 1042.32
$ws.Range("I107").Value = 1121.9474
$ws.Range("J107").Value = 790.1667
$ws.Range("K107").Value = 1121.9474
$ws.Range("L107").Value = 790.1667
$ws.Range("M107").Value = 798.0526
$ws.Range("N107").Value = -4630.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5010000
$ws.Range("J2").Value = 20000
$ws.Range("L2").Value = 20000
$ws.Range("N2").Value = -20224
$ws.Range("H61").Value = 4637980
$ws.Range("I61").Value = 5299493.5
$ws.Range("J61").Value = 7387
$ws.Range("K61").Value = 5299493.5
$ws.Range("L61").Value = 7387
$ws.Range("M61").Value = -5299291.5
$ws.Range("N61").Value = -7791
$ws.Range("H68").Value = 4000.8
$ws.Range("I68").Value = 3751
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 3751
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -3002
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 4000.8
$ws.Range("I71").Value = 3751
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 18755
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -15011
$ws.Range("N71").Value = -32488
$ws.Range("H113").Value = 4637980
$ws.Range("I113").Value = 5299493.5
$ws.Range("J113").Value = 7387
$ws.Range("K113").Value = 5299493.5
$ws.Range("L113").Value = 7387
$ws.Range("M113").Value = -5297323.5
$ws.Range("N113").Value = -11727
$ws.Range("H132").Value = 5568.1577
$ws.Range("I132").Value = 6049.375
$ws.Range("J132").Value = 3001.6667
$ws.Range("K132").Value = 18148.125
$ws.Range("L132").Value = 9005.000100000001
$ws.Range("M132").Value = -15618.125
$ws.Range("N132").Value = -14065.0001
$ws.Range("H136").Value = 36302.03
$ws.Range("I136").Value = 61994.06
$ws.Range("J136").Value = 7184.4
$ws.Range("K136").Value = 185982.18
$ws.Range("L136").Value = 21553.2
$ws.Range("M136").Value = -183432.18
$ws.Range("N136").Value = -26653.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 19499.75
$ws.Range("J25").Value = 19999.666
$ws.Range("L25").Value = 19999.666
$ws.Range("N25").Value = -20585.666
$ws.Range("H113").Value = 1109.7222
$ws.Range("I113").Value = 519
$ws.Range("K113").Value = 1557
$ws.Range("M113").Value = 613
$ws.Range("H132").Value = 22751448
$ws.Range("I132").Value = 31254610
$ws.Range("K132").Value = 93763830
$ws.Range("M132").Value = -93761300
$ws.Range("H136").Value = 2008.4894
$ws.Range("I136").Value = 2024.8572
$ws.Range("J136").Value = 1871
$ws.Range("K136").Value = 6074.571599999999
$ws.Range("L136").Value = 5613
$ws.Range("M136").Value = -3524.571599999999
